# Update specific imputed values in the KNN result data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value  = -13.01
$ws.Range("C4").Value  = -13.318
$ws.Range("C7").Value  = -13.43
$ws.Range("C8").Value  = -12.672
$ws.Range("A11").Value = -21.803
$ws.Range("A12").Value = -21.694
$ws.Range("C12").Value = -13.232
$ws.Range("C14").Value = -12.048
$ws.Range("A15").Value = -21.098
$ws.Range("C22").Value = -13.318
